$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric must be forced to Text format
# first, to preserve the original Text cell type (matches source diff).
$textCells = @("D5","D6","D8","D9","D10","D14","D15","D16","D18","D20","D23","D25","D26","D27","D31","D35","D36","D37","D39","D40","D43","D44","D47","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the crypto price refresh
$ws.Range("D2").Value = "27.571.06"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "1.665.29"
$ws.Range("E3").Value = "  -3.56%  "
$ws.Range("D5").Value = "215.39"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "23.52"
$ws.Range("E8").Value = "  -2.47%  "
$ws.Range("D9").Value = "0.263"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "0.0621"
$ws.Range("E10").Value = "  -1.89%  "
$ws.Range("D12").Value = "1.902.74"
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("D13").Value = "1.673.17"
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").Value = "0.558"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "66.03"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "27.589.89"
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("D18").Value = "241.91"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").Value = "  -3.63%  "
$ws.Range("D20").Value = "7.54"
$ws.Range("E20").Value = "  -4.23%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  -3.32%  "
$ws.Range("D23").Value = "9.27"
$ws.Range("E23").Value = "  -4.89%  "
$ws.Range("D25").Value = "146.12"
$ws.Range("E25").Value = "  -1.96%  "
$ws.Range("D26").Value = "7.17"
$ws.Range("E26").Value = "  -4.76%  "
$ws.Range("D27").Value = "16.34"
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("E30").Value = "  +3.58%  "
$ws.Range("D31").Value = "0.0503"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("D33").Value = "1.477.30"
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  -5.84%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "2.38"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "0.931"
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "0.571"
$ws.Range("E39").Value = "  -6.11%  "
$ws.Range("D40").Value = "69.40"
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("E41").Value = "  -5.53%  "
$ws.Range("D43").Value = "5.39"
$ws.Range("E43").Value = "  -7.56%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.21"
$ws.Range("E44").Value = "  -4.02%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.809.75"
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").Value = "1.70"
$ws.Range("E47").Value = "  -3.40%  "
$ws.Range("D48").Value = "89.24"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  -3.78%  "
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("D51").Value = "7.86"
$ws.Range("E51").Value = "  -3.78%  "
